$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.984.03"
$ws.Range("E2").Value = "  -4.01%  "
$ws.Range("D3").Value = "3.139.12"
$ws.Range("E3").Value = "  -3.52%  "
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").Value = "'605.43"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("D6").Value = "'146.54"
$ws.Range("E6").Value = "  -7.01%  "
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").Value = "3.134.01"
$ws.Range("E8").Value = "  -3.69%  "
$ws.Range("E10").Value = "  -7.57%  "
$ws.Range("D11").Value = "'5.44"
$ws.Range("E11").Value = "  -6.88%  "
$ws.Range("D12").Value = "'0.471"
$ws.Range("E12").Value = "  -6.08%  "
$ws.Range("D13").Value = "'0.0000249"
$ws.Range("E13").Value = "  -8.23%  "
$ws.Range("D14").Value = "'35.26"
$ws.Range("E14").Value = "  -10.13%  "
$ws.Range("D15").Value = "3.661.41"
$ws.Range("E15").Value = "  -3.25%  "
$ws.Range("D16").Value = "64.065.04"
$ws.Range("E16").Value = "  -3.86%  "
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("D18").Value = "3.141.82"
$ws.Range("E18").Value = "  -3.44%  "
$ws.Range("D19").Value = "'6.86"
$ws.Range("E19").Value = "  -8.14%  "
$ws.Range("D20").Value = "'475.87"
$ws.Range("E20").Value = "  -6.36%  "
$ws.Range("D21").Value = "'14.73"
$ws.Range("E21").Value = "  -4.81%  "
$ws.Range("D22").Value = "'0.706"
$ws.Range("E22").Value = "  -6.34%  "
$ws.Range("D23").Value = "'7.75"
$ws.Range("E23").Value = "  -4.74%  "
$ws.Range("D24").Value = "'13.56"
$ws.Range("E24").Value = "  -8.34%  "
$ws.Range("D25").Value = "'82.93"
$ws.Range("E25").Value = "  -4.75%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E27").Value = "  -5.22%  "
$ws.Range("E28").Value = "  -8.07%  "
$ws.Range("E29").Value = "  -10.16%  "
$ws.Range("D30").Value = "'6.75"
$ws.Range("E30").Value = "  -2.31%  "
$ws.Range("E31").Value = "  -16.41%  "
$ws.Range("E32").Value = "  -6.05%  "
$ws.Range("E33").Value = "  +0.20%  "
$ws.Range("D34").Value = "'26.04"
$ws.Range("E34").Value = "  -7.41%  "
$ws.Range("E35").Value = "  -4.99%  "
$ws.Range("D36").Value = "'5.94"
$ws.Range("E36").Value = "  -8.34%  "
$ws.Range("D37").Value = "'53.64"
$ws.Range("E37").Value = "  -3.67%  "
$ws.Range("D38").Value = "0.0₃0732"
$ws.Range("E38").Value = "  -6.87%  "
$ws.Range("D39").Value = "'461.17"
$ws.Range("E39").Value = "  -6.71%  "
$ws.Range("D40").Value = "'2.93"
$ws.Range("E40").Value = "  -14.34%  "
$ws.Range("D41").Value = "'0.0393"
$ws.Range("E41").Value = "  -7.90%  "
$ws.Range("B42").Value = "Cosmos"
$ws.Range("C42").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D42").Value = "'8.39"
$ws.Range("E42").Value = "  -5.11%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "'0.118"
$ws.Range("E43").Value = "  -8.64%  "
$ws.Range("D44").Value = "2.840.19"
$ws.Range("E44").Value = "  -4.99%  "
$ws.Range("D45").Value = "'0.264"
$ws.Range("E45").Value = "  -10.31%  "
$ws.Range("E46").Value = "  -10.27%  "
$ws.Range("D47").Value = "'26.33"
$ws.Range("E47").Value = "  -8.85%  "
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("E49").Value = "  -7.97%  "
$ws.Range("E50").Value = "  -5.15%  "
$ws.Range("D51").Value = "'119.16"
$ws.Range("E51").Value = "  -1.42%  "
